$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(182, 1).Value = 181
$ws.Cells.Item(182, 2).Value = 1
$ws.Cells.Item(182, 3).Value = "2024-06-18 19:10:29"
$ws.Cells.Item(182, 4).Value = 200
$ws.Cells.Item(182, 5).Value = 14

$ws.Cells.Item(183, 1).Value = 182
$ws.Cells.Item(183, 2).Value = 2
$ws.Cells.Item(183, 3).Value = "2024-06-18 19:10:30"
$ws.Cells.Item(183, 4).Value = 200
$ws.Cells.Item(183, 5).Value = 2
